$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "I have faced significant barriers in finding relevant opportunities to apply my programming skills"
$ws.Range("B2").Value = 0.2991889419943792

$ws.Range("A3").Value = "How would you describe your current perception of learning computer programming?_Limited awareness: I have limited knowledge or understanding of programming and its benefits;"
$ws.Range("B3").Value = 0.06826131483818645

$ws.Range("A4").Value = "How would you rate the level of support you receive from your parent/guardian/spouse in learning computer programming?"
$ws.Range("B4").Value = 0.05472762994605947

$ws.Range("A5").Value = "Have you found online coding platforms or communities specifically tailored to the needs of young Africans?"
$ws.Range("B5").Value = 0.03627575334464839

$ws.Range("A6").Value = "I am aware of the various career opportunities available for individuals with computer programming skills"
$ws.Range("B6").Value = 0.03605868920991648

$ws.Range("A7").Value = "How would you describe your current perception of learning computer programming?_Exciting and valuable: Learning programming is an exciting opportunity that provides valuable skills for the future;Challenging but worthwhile: Learning programming can be challenging, but the benefits and rewards make it worth the effort;"
$ws.Range("B7").Value = 0.02892735437053113

$ws.Range("A8").Value = "On a scale of 1 to 10, how interested are you in learning computer programming or improving on your skills?_10"
$ws.Range("B8").Value = 0.02671757396340805

$ws.Range("A9").Value = "Programming is primarily for individuals with strong mathematical or technical backgrounds"
$ws.Range("B9").Value = 0.02214951105968853

$ws.Range("A10").Value = "Gender"
$ws.Range("B10").Value = 0.02188639950398212

$ws.Range("A11").Value = "Country of residence_Nigeria"
$ws.Range("B11").Value = 0.02162327023633552

$ws.Range("A12").Value = "Employment status"
$ws.Range("B12").Value = 0.020813054788948

$ws.Range("A13").Value = "I have encountered societal beliefs that discourage my pursuit of computer programming"
$ws.Range("B13").Value = 0.01995918008489123

$ws.Range("A14").Value = "It is important for me to see more representation of my cultural or societal background in the programming industry"
$ws.Range("B14").Value = 0.01824195851240663

$ws.Range("A15").Value = "My immediate family is supportive of my decision to pursue programming"
$ws.Range("B15").Value = 0.01811134141549349

$ws.Range("A16").Value = "I have experienced societal pressure to pursue traditional career paths instead of programming"
$ws.Range("B16").Value = 0.01777049485263905

$ws.Range("A17").Value = "I have faced gender-related barriers or stereotypes that discourage my involvement in programming"
$ws.Range("B17").Value = 0.01674723479807783

$ws.Range("A18").Value = "What career paths are you interested in pursuing with your computer programming skills?_Not sure"
$ws.Range("B18").Value = 0.01665973270505237

$ws.Range("A19").Value = "It is important for me to have access to mentors or role models in the field of computer programming"
$ws.Range("B19").Value = 0.01564709726321841

$ws.Range("A20").Value = "Age"
$ws.Range("B20").Value = 0.01537291076540824

$ws.Range("A21").Value = "Do you have access to a computer or laptop at your home or within your community?"
$ws.Range("B21").Value = 0.01364846303272846
